$wb = $excel.ActiveWorkbook

# ALC row 38
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 1230.7778
$ws.Range("J38").Value = 1918.7778
$ws.Range("L38").Value = 5756.3334
$ws.Range("N38").Value = -6500.3334

# ALC row 58
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H58").Value = 2466
$ws.Range("J58").Value = 4179.6665
$ws.Range("L58").Value = 12538.9995
$ws.Range("N58").Value = -12838.9995

# ALC row 69
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 5832.857
$ws.Range("J69").Value = 5707.5
$ws.Range("L69").Value = 17122.5
$ws.Range("N69").Value = -18870.5

# ALC row 72
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 5832.857
$ws.Range("J72").Value = 5707.5
$ws.Range("L72").Value = 51367.5
$ws.Range("N72").Value = -60103.5

# ALC row 80
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 629.5294
$ws.Range("I80").Value = 1176.6666
$ws.Range("J80").Value = 512.2857
$ws.Range("K80").Value = 3529.9998
$ws.Range("L80").Value = 1536.8571
$ws.Range("M80").Value = -2531.9998
$ws.Range("N80").Value = -3532.8571

# ALC row 83
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H83").Value = 629.5294
$ws.Range("I83").Value = 1176.6666
$ws.Range("J83").Value = 512.2857
$ws.Range("K83").Value = 10589.9994
$ws.Range("L83").Value = 4610.571300000001
$ws.Range("M83").Value = -5597.999400000001
$ws.Range("N83").Value = -14594.5713

# ALC row 94
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H94").Value = 2501.25
$ws.Range("I94").Value = 2501.25
$ws.Range("K94").Value = 2501.25
$ws.Range("M94").Value = -2050.25

# ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 1222.6833
$ws.Range("I129").Value = 922
$ws.Range("K129").Value = 2766
$ws.Range("M129").Value = 2234

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9386.879000000001
$ws.Range("I32").Value = 8240.063
$ws.Range("J32").Value = 16936.75
$ws.Range("K32").Value = 8240.063
$ws.Range("L32").Value = 16936.75
$ws.Range("M32").Value = -7953.063
$ws.Range("N32").Value = -17510.75

# CRP row 4
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 9572.5
$ws.Range("J4").Value = 10082.857
$ws.Range("L4").Value = 10082.857
$ws.Range("N4").Value = -10306.857

# CRP row 99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 5884889.5
$ws.Range("I99").Value = 2428.7273
$ws.Range("J99").Value = 16669400
$ws.Range("K99").Value = 2428.7273
$ws.Range("L99").Value = 16669400
$ws.Range("M99").Value = -930.7273
$ws.Range("N99").Value = -16672396

# CRP row 126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 5884889.5
$ws.Range("I126").Value = 2428.7273
$ws.Range("J126").Value = 16669400
$ws.Range("K126").Value = 7286.1819
$ws.Range("L126").Value = 50008200
$ws.Range("M126").Value = -4816.1819
$ws.Range("N126").Value = -50013140

# CUL row 4
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 6511

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1043.3448
$ws.Range("I131").Value = 612.8570999999999
$ws.Range("J131").Value = 1180.3182
$ws.Range("K131").Value = 1838.5713
$ws.Range("L131").Value = 3540.9546
$ws.Range("M131").Value = 3201.4287
$ws.Range("N131").Value = -13620.9546

# GSM row 5
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()

# GSM row 69
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 20000
$ws.Range("J69").Value = 20000
$ws.Range("L69").Value = 20000
$ws.Range("N69").Value = -21498

# GSM row 72
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H72").Value = 20000
$ws.Range("J72").Value = 20000
$ws.Range("L72").Value = 60000
$ws.Range("N72").Value = -67488

# GSM row 80
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2311.9443
$ws.Range("I80").Value = 2300.3572
$ws.Range("J80").Value = 2352.5
$ws.Range("K80").Value = 2300.3572
$ws.Range("L80").Value = 2352.5
$ws.Range("M80").Value = -1302.3572
$ws.Range("N80").Value = -4348.5

# GSM row 83
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 2311.9443
$ws.Range("I83").Value = 2300.3572
$ws.Range("J83").Value = 2352.5
$ws.Range("K83").Value = 11501.786
$ws.Range("L83").Value = 11762.5
$ws.Range("M83").Value = -6509.786
$ws.Range("N83").Value = -21746.5

# GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 2318801.2
$ws.Range("I132").Value = 5955699.5
$ws.Range("J132").Value = 4411.4546
$ws.Range("K132").Value = 17867098.5
$ws.Range("L132").Value = 13234.3638
$ws.Range("M132").Value = -17864568.5
$ws.Range("N132").Value = -18294.3638

# LTW row 2
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 2312920.5
$ws.Range("J2").Value = 2312920.5
$ws.Range("L2").Value = 2312920.5
$ws.Range("N2").Value = -2313144.5

# LTW row 7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3005.6667
$ws.Range("I7").Value = 2763.5
$ws.Range("J7").Value = 3490
$ws.Range("K7").Value = 2763.5
$ws.Range("L7").Value = 3490
$ws.Range("M7").Value = -2651.5
$ws.Range("N7").Value = -3714

# LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3460.3157
$ws.Range("I40").Value = 2897.0667
$ws.Range("K40").Value = 2897.0667
$ws.Range("M40").Value = -2761.0667

# LTW row 126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 3005.6667
$ws.Range("I126").Value = 2763.5
$ws.Range("J126").Value = 3490
$ws.Range("K126").Value = 8290.5
$ws.Range("L126").Value = 10470
$ws.Range("M126").Value = -5820.5
$ws.Range("N126").Value = -15410

# LTW row 132
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 50877.547
$ws.Range("I132").Value = 96409.82000000001
$ws.Range("J132").Value = 5345.273
$ws.Range("K132").Value = 289229.46
$ws.Range("L132").Value = 16035.819
$ws.Range("M132").Value = -286699.46
$ws.Range("N132").Value = -21095.819

# WVR row 62
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4662.5713
$ws.Range("I62").Value = 3750
$ws.Range("J62").Value = 4814.6665
$ws.Range("K62").Value = 3750
$ws.Range("L62").Value = 4814.6665
$ws.Range("M62").Value = -3126
$ws.Range("N62").Value = -6062.6665

# WVR row 65
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H65").Value = 4662.5713
$ws.Range("I65").Value = 3750
$ws.Range("J65").Value = 4814.6665
$ws.Range("K65").Value = 18750
$ws.Range("L65").Value = 24073.3325
$ws.Range("M65").Value = -15630
$ws.Range("N65").Value = -30313.3325

# WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 1851.3
$ws.Range("I81").Value = 1846.2
$ws.Range("J81").Value = 1866.6
$ws.Range("K81").Value = 3692.4
$ws.Range("L81").Value = 3733.2
$ws.Range("M81").Value = -2631.4
$ws.Range("N81").Value = -5855.2

# WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 1851.3
$ws.Range("I84").Value = 1846.2
$ws.Range("J84").Value = 1866.6
$ws.Range("K84").Value = 18462
$ws.Range("L84").Value = 18666
$ws.Range("M84").Value = -13158
$ws.Range("N84").Value = -29274

# WVR row 126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 2421.884
$ws.Range("I126").Value = 2352.8794
$ws.Range("K126").Value = 7058.638199999999
$ws.Range("M126").Value = -4588.638199999999

# WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1944.9025
$ws.Range("I132").Value = 2123.4119
$ws.Range("J132").Value = 1818.4584
$ws.Range("K132").Value = 6370.2357
$ws.Range("L132").Value = 5455.3752
$ws.Range("M132").Value = -3840.2357
$ws.Range("N132").Value = -10515.3752
